# Update the cryptocurrency price/volume snapshot for this run.
# Column D values are stored as text (some prices contain two "."
# thousand-separators, e.g. "27.191.01", which is not a valid number),
# so a leading apostrophe is used to force text entry, and the style
# is reset afterward so the quote-prefix marker is not retained.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.191.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "'1.904.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'306.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.5240"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("D8").Value = "'0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.07255"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "'0.9035"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "'0.08510"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.37%  "
$ws.Range("D13").Value = "'1.930.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").Value = "'95.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "'5.297"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'0.000008636"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'27.231.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'5.072"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'2.153.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "'6.435"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.298"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'147.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.750"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").Value = "'115.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.817"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.916"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'0.09280"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").Value = "'0.8054"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "'0.05051"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "'1.239"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'3.449"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.86%  "
$ws.Range("D37").Value = "'2.960"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'2.625"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'0.5710"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'9.192"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "'6.650"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'116.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "'0.1520"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'0.4870"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").Value = "'1.615"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "'37.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'64.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
